$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns stay text (they contain values like "30.135.29"
# and "42.69" which Excel would otherwise silently coerce into numbers/dates).
$ws.Columns("D:E").NumberFormat = "@"

$rows = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "30.135.29", "  -0.60%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.860.89", "  -0.44%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "0.9987", "  -0.27%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "241.47", "  +2.79%  "),
    @(6, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "0.9984", "  -0.24%  "),
    @(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.4683", "  -0.39%  "),
    @(8, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "42.69", "  -0.56%  "),
    @(9, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.2857", "  -0.40%  "),
    @(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.06473", "  -1.56%  "),
    @(11, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "20.81", "  -3.76%  "),
    @(12, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.07680", "  -4.15%  "),
    @(13, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.860.85", "  -0.51%  "),
    @(14, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "93.85", "  -3.09%  "),
    @(15, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.067", "  -0.92%  "),
    @(16, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.6797", "  -0.58%  "),
    @(17, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "269.21", "  -0.21%  "),
    @(18, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "30.116.34", "  -0.65%  "),
    @(19, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "13.31", "  -4.99%  "),
    @(20, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000007517", "  -1.43%  "),
    @(21, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.9990", "  -0.19%  "),
    @(22, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.090.41", "  -1.23%  "),
    @(23, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "0.9984", "  -0.22%  "),
    @(24, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "5.161", "  -2.02%  "),
    @(25, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "6.112", "  -1.51%  "),
    @(26, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "9.292", "  -1.11%  "),
    @(27, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "165.84", "  -1.27%  "),
    @(28, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "18.77", "  -0.69%  "),
    @(29, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "1.882", "  -3.36%  "),
    @(30, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.370", "  +0.07%  "),
    @(31, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.09821", "  -0.63%  "),
    @(32, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.500", "  +2.52%  "),
    @(33, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.210", "  -3.53%  "),
    @(34, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "3.989", "  -1.80%  "),
    @(35, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.04666", "  -1.08%  "),
    @(36, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.109", "  -2.37%  "),
    @(37, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.6853", "  -2.06%  "),
    @(38, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.704", "  -0.33%  "),
    @(39, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01826", "  -2.42%  "),
    @(40, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.717", "  +2.26%  "),
    @(41, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "6.319", "  +0.74%  "),
    @(42, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "70.32", "  -2.05%  "),
    @(43, "PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "0.9980", "  -0.20%  "),
    @(44, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.8328", "  -1.12%  "),
    @(45, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "1.885", "  -3.67%  "),
    @(46, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "101.77", "  -0.89%  "),
    @(47, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.4045", "  -2.82%  "),
    @(48, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "9.201", "  +1.08%  "),
    @(49, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "927.70", "  +1.92%  "),
    @(50, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "6.934", "  -1.68%  "),
    @(51, "Elrond", "https://coinranking.com/coin/omwkOTglq+elrond-egld", "34.13", "  -0.84%  ")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
